# Fruta / hortaliza, semanal
# Insert a new weekly price-log entry as row 405, pushing all following rows
# down by one (old row 405 -> 406, ..., old row 473 -> new row 474).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 405.
$ws.Rows.Item(405).Insert()

# Populate the newly inserted row 405 with the new weekly entry.
$ws.Cells.Item(405, 1).Value = 4
$ws.Cells.Item(405, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(405, 3).Value = "Los Lagos"
$ws.Cells.Item(405, 4).Value = 45209
$ws.Cells.Item(405, 5).Value = 10
$ws.Cells.Item(405, 6).Value = "Fruta"
$ws.Cells.Item(405, 7).Value = 100108
$ws.Cells.Item(405, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(405, 9).Value = 100108005
$ws.Cells.Item(405, 10).Value = "Piña"
$ws.Cells.Item(405, 11).Value = "Caramelo"
$ws.Cells.Item(405, 12).Value = "Segunda"
$ws.Cells.Item(405, 13).Value = 150
$ws.Cells.Item(405, 14).Value = 25000
$ws.Cells.Item(405, 15).Value = 25000
$ws.Cells.Item(405, 16).Value = 25000
$ws.Cells.Item(405, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(405, 18).Value = "Ecuador"
$ws.Cells.Item(405, 19).Value = 1786
$ws.Cells.Item(405, 20).Value = 14
